# Update "想去人数" (interest count) / "最低票价" (min price) figures that
# changed between scrapes, across the 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G3").Value = 73
$ws.Range("F4").Value = 961
$ws.Range("F5").Value = 204
$ws.Range("F7").Value = 1051
$ws.Range("F8").Value = 837
$ws.Range("F9").Value = 249
$ws.Range("F12").Value = 840
$ws.Range("F13").Value = 289
$ws.Range("F15").Value = 505
$ws.Range("F16").Value = 1336
$ws.Range("F18").Value = 1204
$ws.Range("F19").Value = 1198
$ws.Range("F20").Value = 2887
$ws.Range("F21").Value = 1439
$ws.Range("F22").Value = 706
$ws.Range("F24").Value = 1278
$ws.Range("F28").Value = 3134
$ws.Range("F29").Value = 611
$ws.Range("F31").Value = 1415

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 368
$ws.Range("F13").Value = 23
$ws.Range("F14").Value = 5

# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 747

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 747
$ws.Range("G4").Value = 73
$ws.Range("F6").Value = 368
$ws.Range("F9").Value = 961
$ws.Range("F10").Value = 204
$ws.Range("F13").Value = 1051
$ws.Range("F14").Value = 837
$ws.Range("F15").Value = 249
$ws.Range("F23").Value = 840
$ws.Range("F24").Value = 289
$ws.Range("F26").Value = 505
$ws.Range("F27").Value = 1336
$ws.Range("F29").Value = 1204
$ws.Range("F30").Value = 1198
$ws.Range("F31").Value = 2887
$ws.Range("F32").Value = 1439
$ws.Range("F33").Value = 706
$ws.Range("F35").Value = 1278
$ws.Range("F38").Value = 23
$ws.Range("F41").Value = 3135
$ws.Range("F42").Value = 611
$ws.Range("F44").Value = 1415
$ws.Range("F45").Value = 5
